$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 3, pushing old rows 3-6 down to 5-8.
$ws.Rows("3:4").Insert()

# Row 3: new weekly entry (Primera, Provincia de Limari)
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = 44462
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = "Fruta"
$ws.Range("G3").Value = 100107
$ws.Range("H3").Value = "Otros"
$ws.Range("I3").Value = 100107002
$ws.Range("J3").Value = "Chirimoya"
$ws.Range("K3").Value = "Cultivar IV Región"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 2900
$ws.Range("O3").Value = 3000
$ws.Range("P3").Value = 2950
$ws.Range("Q3").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 2950
$ws.Range("T3").Value = 1

# Row 4: new weekly entry (Segunda, Provincia de Limari)
$ws.Range("A4").Value = 11
$ws.Range("B4").Value = "Vega Monumental Concepción"
$ws.Range("C4").Value = "Bíobío"
$ws.Range("D4").Value = 44462
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100107
$ws.Range("H4").Value = "Otros"
$ws.Range("I4").Value = 100107002
$ws.Range("J4").Value = "Chirimoya"
$ws.Range("K4").Value = "Cultivar IV Región"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 2600
$ws.Range("O4").Value = 2600
$ws.Range("P4").Value = 2600
$ws.Range("Q4").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 2600
$ws.Range("T4").Value = 1
